$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 21-22, shifting the existing data (old rows 21-33)
# down to rows 23-35.
$ws.Rows("21:22").Insert()

# Row 21: new "Especial" quality record
$ws.Range("A21").Value = 8
$ws.Range("B21").Value = "Terminal La Palmera de La Serena"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44669
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100104
$ws.Range("H21").Value = "Frutos de pepita"
$ws.Range("I21").Value = 100104003
$ws.Range("J21").Value = "Membrillo"
$ws.Range("K21").Value = "Champion"
$ws.Range("L21").Value = "Especial"
$ws.Range("M21").Value = 16
$ws.Range("N21").Value = 330000
$ws.Range("O21").Value = 340000
$ws.Range("P21").Value = 335000
$ws.Range("Q21").Value = "$/bins (450 kilos)"
$ws.Range("R21").Value = "Región de O'Higgins"
$ws.Range("S21").Value = 744
$ws.Range("T21").Value = 450

# Row 22: new "Primera" quality record
$ws.Range("A22").Value = 8
$ws.Range("B22").Value = "Terminal La Palmera de La Serena"
$ws.Range("C22").Value = "Coquimbo"
$ws.Range("D22").Value = 44669
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = "Fruta"
$ws.Range("G22").Value = 100104
$ws.Range("H22").Value = "Frutos de pepita"
$ws.Range("I22").Value = 100104003
$ws.Range("J22").Value = "Membrillo"
$ws.Range("K22").Value = "Champion"
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 16
$ws.Range("N22").Value = 290000
$ws.Range("O22").Value = 300000
$ws.Range("P22").Value = 295000
$ws.Range("Q22").Value = "$/bins (450 kilos)"
$ws.Range("R22").Value = "Región de O'Higgins"
$ws.Range("S22").Value = 656
$ws.Range("T22").Value = 450
